# Weekly update: insert two new daily price records (2021-12-23 and
# 2021-12-24) into the "Poroto granado" sheet, in date order, shifting the
# subsequent rows down (matching how the source dataset is appended to
# weekly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 31 (2021-12-24, value 44554) -----------------------
$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 7
$ws.Range("B31").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C31").Value = "Ñuble"
$ws.Range("D31").Value = 44554
$ws.Range("E31").Value = 16
$ws.Range("F31").Value = 100112030
$ws.Range("G31").Value = "Poroto granado"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 60
$ws.Range("K31").Value = 37000
$ws.Range("L31").Value = 38000
$ws.Range("M31").Value = 37500
$ws.Range("N31").Value = "$/saco 25 kilos"
$ws.Range("O31").Value = "Región del Maule"
$ws.Range("P31").Value = 1500
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"

# --- Insert new row at 45 (2021-12-23, value 44553) -----------------------
$ws.Rows("45:45").Insert()

$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 44553
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112030
$ws.Range("G45").Value = "Poroto granado"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 60
$ws.Range("K45").Value = 35000
$ws.Range("L45").Value = 36000
$ws.Range("M45").Value = 35500
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 1420
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
